$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Job column (C2) for row 2, matching row 3's value
$ws.Range("C2").Value = "IT"

# Update the Domeniu column (D2) with the new value
$ws.Range("D2").Value = "IT-Softwardfghje"

# Update the active selection to D2, as recorded in the saved view state
$ws.Range("D2").Select()
